$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'331.55"
$ws.Range("E2").Value = "'0.21%"
$ws.Range("D3").Value = "'41.82"
$ws.Range("E3").Value = "'1.29%"
$ws.Range("D4").Value = "'5.673"
$ws.Range("E4").Value = "'-0.71%"
$ws.Range("D5").Value = "'0.08383"
$ws.Range("E5").Value = "'3.42%"
$ws.Range("D6").Value = "'8.797"
$ws.Range("E6").Value = "'1.40%"
$ws.Range("D7").Value = "'2.018"
$ws.Range("E7").Value = "'2.50%"
$ws.Range("D8").Value = "'4.529"
$ws.Range("E8").Value = "'0.73%"
$ws.Range("D9").Value = "'2.891"
$ws.Range("E9").Value = "'-3.61%"
$ws.Range("D10").Value = "'0.9272"
$ws.Range("E10").Value = "'-0.08%"
$ws.Range("D11").Value = "'0.1295"
$ws.Range("E11").Value = "'1.41%"
$ws.Range("D12").Value = "'0.1970"
$ws.Range("E12").Value = "'0.49%"
$ws.Range("D13").Value = "'0.09383"
$ws.Range("E13").Value = "'2.13%"
$ws.Range("D14").Value = "'0.03891"
$ws.Range("E14").Value = "'4.06%"
$ws.Range("D15").Value = "'0.1061"
$ws.Range("E15").Value = "'0.90%"
$ws.Range("D16").Value = "'0.001295"
$ws.Range("E16").Value = "'0.31%"
$ws.Range("D17").Value = "'0.006111"
$ws.Range("E17").Value = "'-3.46%"
$ws.Range("E18").Value = "'1.93%"
$ws.Range("D20").Value = "'8.156"
$ws.Range("E20").Value = "'-7.29%"
$ws.Range("D21").Value = "'0.1373"
$ws.Range("E21").Value = "'0.43%"
$ws.Range("E22").Value = "'0.31%"
$ws.Range("D23").Value = "'0.04415"
$ws.Range("E23").Value = "'-0.17%"
$ws.Range("E24").Value = "'0.30%"
$ws.Range("D25").Value = "'0.004405"
$ws.Range("E25").Value = "'-0.11%"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'-3.04%"
$ws.Range("D39").Value = "'0.02798"
$ws.Range("E39").Value = "'0.23%"
$ws.Range("D40").Value = "'0.05531"
$ws.Range("E40").Value = "'-0.55%"
$ws.Range("D41").Value = "'0.007807"
$ws.Range("E41").Value = "'3.36%"
$ws.Range("D42").Value = "'0.1435"
$ws.Range("E42").Value = "'0.99%"
$ws.Range("D43").Value = "'0.009325"
$ws.Range("D44").Value = "'0.002160"
$ws.Range("E44").Value = "'2.65%"
$ws.Range("D45").Value = "'0.01091"
$ws.Range("E45").Value = "'-7.51%"
$ws.Range("D46").Value = "'0.00007017"
$ws.Range("E46").Value = "'3.59%"
$ws.Range("E47").Value = "'0.39%"
$ws.Range("E48").Value = "'14.87%"
$ws.Range("E49").Value = "'0.32%"
$ws.Range("E50").Value = "'0.39%"
$ws.Range("E51").Value = "'0.39%"
